# The "two step login" bullet originally read:
#   --> for practicing work with more AuthenticationProviders.
# The author re-typed the trailing plural "s" (splitting it into separate
# keystrokes/runs: "AuthenticationProvider" + "r" + "s."), which produced a
# duplicated "r" in the final text. Reproduce that exact end state while
# keeping the run's original character formatting intact.

$d = $word.ActiveDocument

$found = $d.Content.Find.Execute(
    "AuthenticationProviders.",  # old text
    $true,                       # MatchCase
    $false,                      # MatchWholeWord
    $false,                      # MatchWildcards
    $false,                      # MatchSoundsLike
    $false,                      # MatchAllWordForms
    $true,                       # Forward
    1,                           # Wrap (wdFindContinue)
    $false,                      # Format
    "AuthenticationProviderrs.", # new text
    2                            # Replace (wdReplaceAll)
)

Write-Host "Replaced: $found"
